$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "13Ballard" entry (Cell Number 13) and move the "28Ballard" row
# (previously row 5) up to row 2, ahead of "18Ballard" and "20Ballard".
# Resulting order: 28Ballard, 18Ballard, 20Ballard.

$ws.Range("A2").Value = "28Ballard"
$ws.Range("B2").Value = 28

$ws.Range("A3").Value = "18Ballard"
$ws.Range("B3").Value = 18

$ws.Range("A4").Value = "20Ballard"
$ws.Range("B4").Value = 20

# The old row 5 (duplicate of the relocated 28Ballard data) is no longer needed.
$ws.Rows.Item(5).Delete()

# Match the saved selection state from the edit.
$ws.Range("A3:XFD3").Select()
